# This workbook contains a single data table (Hortaliza / Sandia prices at
# "Vega Modelo de Temuco"). The edit adds 3 new weekly price observations,
# inserted as new rows right before the existing row 297. Inserting the
# rows shifts all the subsequent rows down by 3 (carrying their formatting
# with them), which matches the rest of the diff (old row N becomes new
# row N+3 for N from 297 to 393), and grows the used range from
# A1:R393 to A1:R396.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 297; existing rows 297-393 shift to 300-396.
$ws.Rows("297:299").Insert()

# --- New row 297 ---
$ws.Cells.Item(297, 1).Value = 10
$ws.Cells.Item(297, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(297, 3).Value = "La Araucanía"
$ws.Cells.Item(297, 4).Value2 = 44559
$ws.Cells.Item(297, 5).Value = 9
$ws.Cells.Item(297, 6).Value = 100112028
$ws.Cells.Item(297, 7).Value = "Sandia"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Extra"
$ws.Cells.Item(297, 10).Value = 950
$ws.Cells.Item(297, 11).Value = 3500
$ws.Cells.Item(297, 12).Value = 3500
$ws.Cells.Item(297, 13).Value = 3500
$ws.Cells.Item(297, 14).Value = "$/unidad"
$ws.Cells.Item(297, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(297, 16).Value = 3500
$ws.Cells.Item(297, 17).Value = 1
$ws.Cells.Item(297, 18).Value = "Hortaliza"

# --- New row 298 ---
$ws.Cells.Item(298, 1).Value = 10
$ws.Cells.Item(298, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(298, 3).Value = "La Araucanía"
$ws.Cells.Item(298, 4).Value2 = 44559
$ws.Cells.Item(298, 5).Value = 9
$ws.Cells.Item(298, 6).Value = 100112028
$ws.Cells.Item(298, 7).Value = "Sandia"
$ws.Cells.Item(298, 8).Value = "Sin especificar"
$ws.Cells.Item(298, 9).Value = "Primera"
$ws.Cells.Item(298, 10).Value = 2500
$ws.Cells.Item(298, 11).Value = 3000
$ws.Cells.Item(298, 12).Value = 3000
$ws.Cells.Item(298, 13).Value = 3000
$ws.Cells.Item(298, 14).Value = "$/unidad"
$ws.Cells.Item(298, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(298, 16).Value = 3000
$ws.Cells.Item(298, 17).Value = 1
$ws.Cells.Item(298, 18).Value = "Hortaliza"

# --- New row 299 ---
$ws.Cells.Item(299, 1).Value = 10
$ws.Cells.Item(299, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(299, 3).Value = "La Araucanía"
$ws.Cells.Item(299, 4).Value2 = 44559
$ws.Cells.Item(299, 5).Value = 9
$ws.Cells.Item(299, 6).Value = 100112028
$ws.Cells.Item(299, 7).Value = "Sandia"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Segunda"
$ws.Cells.Item(299, 10).Value = 3150
$ws.Cells.Item(299, 11).Value = 2500
$ws.Cells.Item(299, 12).Value = 2500
$ws.Cells.Item(299, 13).Value = 2500
$ws.Cells.Item(299, 14).Value = "$/unidad"
$ws.Cells.Item(299, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(299, 16).Value = 2500
$ws.Cells.Item(299, 17).Value = 1
$ws.Cells.Item(299, 18).Value = "Hortaliza"

Write-Host "Inserted 3 new rows (297-299) and populated them."
